$d = $word.ActiveDocument

$pairs = @(
  @("232÷8=", "414÷2="),
  @("164÷6=", "488÷4="),
  @("420÷2=", "519÷3="),
  @("928÷5=", "881÷3="),
  @("137÷9=", "250÷6="),
  @("338÷4=", "612÷2="),
  @("176÷8=", "868÷9="),
  @("767÷2=", "810÷7="),
  @("387÷6=", "134÷3="),
  @("206÷6=", "526÷6="),
  @("952÷7=", "665÷3="),
  @("947÷2=", "642÷3="),
  @("822÷7=", "708÷3="),
  @("937÷5=", "498÷6="),
  @("966÷7=", "822÷4="),
  @("644÷7=", "955÷2="),
  @("660÷4=", "864÷8="),
  @("801÷8=", "129÷2="),
  @("295÷8=", "656÷5="),
  @("997÷9=", "736÷4="),
  @("922÷4=", "782÷6="),
  @("702÷9=", "132÷5="),
  @("724÷6=", "334÷9="),
  @("279÷4=", "337÷2="),
  @("367÷5=", "889÷2=")
)

foreach ($pair in $pairs) {
  $old = $pair[0]
  $new = $pair[1]
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
